$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.2708946666666667
$ws.Cells.Item(2, 8).Value = 0.812684
$ws.Cells.Item(2, 9).Value = 0.1616296696421007
$ws.Cells.Item(2, 10).Value = 0.1616296696421007
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.042868
$ws.Cells.Item(2, 14).Value = 0.128604
$ws.Cells.Item(2, 15).Value = 0.03014606792405771
$ws.Cells.Item(2, 16).Value = 0.03014606792405771
$ws.Cells.Item(2, 17).Value = 0.01161271257066667
$ws.Cells.Item(2, 18).Value = 0.104514413136
$ws.Cells.Item(2, 19).Value = 0.004872498999573775
$ws.Cells.Item(2, 20).Value = 0.004872498999573776
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.2708946666666667
$ws.Cells.Item(3, 8).Value = 0.812684
$ws.Cells.Item(3, 9).Value = 0.1616296696421007
$ws.Cells.Item(3, 10).Value = 0.1616296696421007
$ws.Cells.Item(3, 15).Value = 0.2718481285523376
$ws.Cells.Item(3, 16).Value = 0.2718481285523376
$ws.Cells.Item(3, 17).Value = 0.1047199318897778
$ws.Cells.Item(3, 18).Value = 0.942479387008
$ws.Cells.Item(3, 19).Value = 0.04393872321073763
$ws.Cells.Item(3, 20).Value = 0.04393872321073763
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.2708946666666667
$ws.Cells.Item(4, 8).Value = 0.812684
$ws.Cells.Item(4, 9).Value = 0.1616296696421007
$ws.Cells.Item(4, 10).Value = 0.1616296696421007
$ws.Cells.Item(4, 13).Value = 0.9839956666666666
$ws.Cells.Item(4, 14).Value = 2.951987
$ws.Cells.Item(4, 15).Value = 0.69197537100662
$ws.Cells.Item(4, 16).Value = 0.69197537100662
$ws.Cells.Item(4, 17).Value = 0.2665591781231111
$ws.Cells.Item(4, 18).Value = 2.399032603108
$ws.Cells.Item(4, 19).Value = 0.11184375061627
$ws.Cells.Item(4, 20).Value = 0.11184375061627
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.2708946666666667
$ws.Cells.Item(5, 8).Value = 0.812684
$ws.Cells.Item(5, 9).Value = 0.1616296696421007
$ws.Cells.Item(5, 10).Value = 0.1616296696421007
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.008575333333333332
$ws.Cells.Item(5, 14).Value = 0.025726
$ws.Cells.Item(5, 15).Value = 0.006030432516984765
$ws.Cells.Item(5, 16).Value = 0.006030432516984765
$ws.Cells.Item(5, 17).Value = 0.002323012064888889
$ws.Cells.Item(5, 18).Value = 0.020907108584
$ws.Cells.Item(5, 19).Value = 0.000974696815519229
$ws.Cells.Item(5, 20).Value = 0.0009746968155192291
$ws.Cells.Item(6, 9).Value = 0.6313295261673385
$ws.Cells.Item(6, 10).Value = 0.6313295261673384
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.042868
$ws.Cells.Item(6, 14).Value = 0.128604
$ws.Cells.Item(6, 15).Value = 0.03014606792405771
$ws.Cells.Item(6, 16).Value = 0.03014606792405771
$ws.Cells.Item(6, 17).Value = 0.04535954531733333
$ws.Cells.Item(6, 18).Value = 0.408235907856
$ws.Cells.Item(6, 19).Value = 0.01903210277830376
$ws.Cells.Item(6, 20).Value = 0.01903210277830376
$ws.Cells.Item(7, 9).Value = 0.6313295261673385
$ws.Cells.Item(7, 10).Value = 0.6313295261673384
$ws.Cells.Item(7, 15).Value = 0.2718481285523376
$ws.Cells.Item(7, 16).Value = 0.2718481285523376
$ws.Cells.Item(7, 19).Value = 0.171625750188425
$ws.Cells.Item(7, 20).Value = 0.171625750188425
$ws.Cells.Item(8, 9).Value = 0.6313295261673385
$ws.Cells.Item(8, 10).Value = 0.6313295261673384
$ws.Cells.Item(8, 13).Value = 0.9839956666666666
$ws.Cells.Item(8, 14).Value = 2.951987
$ws.Cells.Item(8, 15).Value = 0.69197537100662
$ws.Cells.Item(8, 16).Value = 0.69197537100662
$ws.Cells.Item(8, 17).Value = 1.041186806807556
$ws.Cells.Item(8, 18).Value = 9.370681261268
$ws.Cells.Item(8, 19).Value = 0.4368644830970776
$ws.Cells.Item(8, 20).Value = 0.4368644830970776
$ws.Cells.Item(9, 9).Value = 0.6313295261673385
$ws.Cells.Item(9, 10).Value = 0.6313295261673384
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.008575333333333332
$ws.Cells.Item(9, 14).Value = 0.025726
$ws.Cells.Item(9, 15).Value = 0.006030432516984765
$ws.Cells.Item(9, 16).Value = 0.006030432516984765
$ws.Cells.Item(9, 17).Value = 0.009073743140444445
$ws.Cells.Item(9, 18).Value = 0.081663688264
$ws.Cells.Item(9, 19).Value = 0.003807190103532102
$ws.Cells.Item(9, 20).Value = 0.003807190103532102
$ws.Cells.Item(10, 9).Value = 0.2070408041905609
$ws.Cells.Item(10, 10).Value = 0.2070408041905609
$ws.Cells.Item(10, 11).Value = 1
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.042868
$ws.Cells.Item(10, 14).Value = 0.128604
$ws.Cells.Item(10, 15).Value = 0.03014606792405771
$ws.Cells.Item(10, 16).Value = 0.03014606792405771
$ws.Cells.Item(10, 17).Value = 0.01487539605066667
$ws.Cells.Item(10, 18).Value = 0.133878564456
$ws.Cells.Item(10, 19).Value = 0.006241466146180181
$ws.Cells.Item(10, 20).Value = 0.006241466146180182
$ws.Cells.Item(11, 9).Value = 0.2070408041905609
$ws.Cells.Item(11, 10).Value = 0.2070408041905609
$ws.Cells.Item(11, 15).Value = 0.2718481285523376
$ws.Cells.Item(11, 16).Value = 0.2718481285523376
$ws.Cells.Item(11, 19).Value = 0.05628365515317495
$ws.Cells.Item(11, 20).Value = 0.05628365515317495
$ws.Cells.Item(12, 9).Value = 0.2070408041905609
$ws.Cells.Item(12, 10).Value = 0.2070408041905609
$ws.Cells.Item(12, 13).Value = 0.9839956666666666
$ws.Cells.Item(12, 14).Value = 2.951987
$ws.Cells.Item(12, 15).Value = 0.69197537100662
$ws.Cells.Item(12, 16).Value = 0.69197537100662
$ws.Cells.Item(12, 17).Value = 0.3414510883131111
$ws.Cells.Item(12, 18).Value = 3.073059794818
$ws.Cells.Item(12, 19).Value = 0.1432671372932723
$ws.Cells.Item(12, 20).Value = 0.1432671372932723
$ws.Cells.Item(13, 9).Value = 0.2070408041905609
$ws.Cells.Item(13, 10).Value = 0.2070408041905609
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.008575333333333332
$ws.Cells.Item(13, 14).Value = 0.025726
$ws.Cells.Item(13, 15).Value = 0.006030432516984765
$ws.Cells.Item(13, 16).Value = 0.006030432516984765
$ws.Cells.Item(13, 17).Value = 0.002975680684888889
$ws.Cells.Item(13, 18).Value = 0.026781126164
$ws.Cells.Item(13, 19).Value = 0.001248545597933434
$ws.Cells.Item(13, 20).Value = 0.001248545597933434
